$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing single data column (A,B) one column to the right,
# so the old "Question Source" header (B1) lands in C1 and the old
# "Logarithms" value (A2) lands in B2 - matching the new layout.
$ws.Columns("A").Insert()

# New column A: "No" / numbering 1..10
$ws.Range("A1").Value = "No"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10

# New column B header
$ws.Range("B1").Value = "Category"

# New column D: "Formula" header + formula text samples
$ws.Range("D1").Value = "Formula"
$ws.Range("D2").Value = "\log_ab"
$ws.Range("D3").Value = "\frac{1}{\log_ab}"

# Column C (old "Question Source" column) gets a text number-format
# and the question-source code values.
$ws.Range("C1:C3").NumberFormat = "@"
$ws.Range("C2").Value = "200604003003"
$ws.Range("C3").Value = "200604003003"

# Column widths: C keeps 14.5 (no longer "best fit"), D is wide for formulas.
$ws.Columns("C").ColumnWidth = 13.666666666666666
$ws.Columns("D").ColumnWidth = 58.5

$ws.Range("C6").Select()
